$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (2-6) got cyclically reshuffled: the row that used to be
# at source row now lives at the target row. Capture the "before" values
# for the columns that vary per-row, then write them back in the new order.

$cols = @("D","M","N","O","P","R","S")
# target row -> source row (which row's original data now appears there)
$map = @{ 2 = 5; 3 = 6; 4 = 3; 5 = 4; 6 = 2 }

$snapshot = @{}
foreach ($row in 2..6) {
    $rowData = @{}
    foreach ($col in $cols) {
        $rowData[$col] = $ws.Range("$col$row").Value2
    }
    $snapshot[$row] = $rowData
}

foreach ($row in 2..6) {
    $srcRow = $map[$row]
    $rowData = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$row").Value = $rowData[$col]
    }
}
